# Estadisticos Matutinos 15 Oct
# Add new "Rescatables" (rescue/makeup exam) records to the Rescatables sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

# Column A - NC (student id numbers)
$ws.Range("A2").Value = 20330051920223
$ws.Range("A3").Value = 20330051920253
$ws.Range("A4").Value = 20330051920298
$ws.Range("A5").Value = 20330051920311

# Column B - Paterno
$ws.Range("B2").Value = "CASTELLANOS"
$ws.Range("B3").Value = "SANCHEZ"
$ws.Range("B4").Value = "LEYVA"
$ws.Range("B5").Value = "ROJAS"

# Column C - Materno
$ws.Range("C2").Value = "TEQUIHUATLE"
$ws.Range("C3").Value = "QUIAHUA"
$ws.Range("C4").Value = "VELAZQUEZ"
$ws.Range("C5").Value = "ROJAS"

# Column D - Nombres
$ws.Range("D2").Value = "JENNIFER"
$ws.Range("D3").Value = "ROSARIO"
$ws.Range("D4").Value = "ELIAN"
$ws.Range("D5").Value = "DULCE MARIA"

# Column E - Nombre_Largo
$ws.Range("E2").Value = "IDENTIFICA MICROORGANISMOS CON BASE EN TÉCNICAS PARASITOLÓGICAS"
$ws.Range("E3").Value = "IDENTIFICA MICROORGANISMOS CON BASE EN TÉCNICAS PARASITOLÓGICAS"
$ws.Range("E4").Value = "IDENTIFICA MICROORGANISMOS CON BASE EN TÉCNICAS BACTERIOLÓGICAS"
$ws.Range("E5").Value = "IDENTIFICA MICROORGANISMOS CON BASE EN TÉCNICAS BACTERIOLÓGICAS"

# Column F - Grupo
$ws.Range("F2").Value = "3ALCM"
$ws.Range("F3").Value = "3ALCM"
$ws.Range("F4").Value = "3BLCM"
$ws.Range("F5").Value = "3BLCM"

# Column G - Reprobadas
$ws.Range("G2").Value = 6
$ws.Range("G3").Value = 6
$ws.Range("G4").Value = 6
$ws.Range("G5").Value = 6
